# Generate Report for Archive
# Update status of the two files that have moved from "Ready for handoff"
# to "In Translation" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: column B = zh-cn status, column C = de-de status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = $newStatus
$overview.Range("C8").Value = $newStatus
$overview.Range("B9").Value = $newStatus
$overview.Range("C9").Value = $newStatus

# --- zh-cn sheet: column B = Status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B8").Value = $newStatus
$zhcn.Range("B9").Value = $newStatus

# --- de-de sheet: column B = Status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B8").Value = $newStatus
$dede.Range("B9").Value = $newStatus
